# Updates the cryptos worksheet Price (D) and Volume(1h) (E) columns
# to reflect the latest GitHub Actions scrape values.
# Numeric-looking Price values are entered with a leading apostrophe
# (forcing text, matching the original text-typed cells) and the
# cell style is reset to Normal afterward so no stray number format
# is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.618.04"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.755.32"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Formula = "'324.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Formula = "'0.4505"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.66%  "
$ws.Range("D8").Formula = "'0.3552"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").Formula = "'0.07455"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").Formula = "'41.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("D11").Formula = "'1.082"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Formula = "'20.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Formula = "'5.981"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").Formula = "'7.143"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "1.749.26"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Formula = "'93.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").Formula = "'0.06457"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Formula = "'17.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").Formula = "'5.747"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").Value = "27.659.48"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Formula = "'2.090"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Formula = "'164.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").Value = "1.954.05"
$ws.Range("D29").Formula = "'2.083"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.07%  "
$ws.Range("D30").Formula = "'125.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Formula = "'1.103"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").Formula = "'0.09196"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.73%  "
$ws.Range("D33").Formula = "'3.659"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").Formula = "'5.491"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("D35").Formula = "'0.02288"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").Formula = "'11.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.16%  "
$ws.Range("D37").Formula = "'0.06030"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Formula = "'0.2081"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("D39").Formula = "'0.6286"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Formula = "'4.931"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").Formula = "'1.183"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Formula = "'1.392"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Formula = "'7.755"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("D44").Formula = "'13.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D45").Formula = "'3.713"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Formula = "'0.5859"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Formula = "'122.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Formula = "'1.936"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Formula = "'0.06888"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").Formula = "'1.129"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").Formula = "'71.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.19%  "